$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Benchmarks")

# I13: literal value -> formula
$ws.Range("I13").Formula = "=1.16/4"

# I14: literal value -> formula
$ws.Range("I14").Formula = "=500/4"

# I16: literal value -> formula, number format changes from "0" to "0.0"
$ws.Range("I16").Formula = "=651/4"
$ws.Range("I16").NumberFormat = "0.0"

# I17: literal value -> formula, number format changes from "0" to "0.0"
$ws.Range("I17").Formula = "=55.791/4"
$ws.Range("I17").NumberFormat = "0.0"

# I18: literal value -> formula, number format changes from "0" to "0.0"
$ws.Range("I18").Formula = "=15.845/4"
$ws.Range("I18").NumberFormat = "0.0"

# I19: was empty -> formula + new K19 formula, number format changes from "0" to "0.0"
$ws.Range("I19").Formula = "=4.617/4"
$ws.Range("I19").NumberFormat = "0.0"
$ws.Range("K19").Formula = "=I19/`$I`$13"

# I20: still empty, but number format changes from "0" to "0.0"
$ws.Range("I20").NumberFormat = "0.0"

# Reset the view: drop the scrolled topLeftCell / prior selection (K38) left over
# from the original file, back to the sheet's home cell.
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("A1").Select()
